$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 6: re-style the sources-of-finance table with the built-in
#    "Medium Style 2 - Accent 1" table style (was the deck's custom
#    "Table_0" style).
# ---------------------------------------------------------------------------
$tableShape = $p.Slides.Item(6).Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B6CF3B1C-2BCE-4C7A-AA67-467637D311A8}", $true)

# ---------------------------------------------------------------------------
# 2. Swap the presentation's theme colours from the custom "Integral"
#    palette over to the standard Office theme palette (Design > Colors >
#    Office). The theme's colour scheme is shared by every slide, so
#    updating it once via any slide's ThemeColorScheme repaints the whole
#    deck.
# ---------------------------------------------------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# Index : role      : target RGB (VBA RGB() little-endian int; R + G*256 + B*65536)
$officeRgb = @(
    0,          # 1  dk1      #000000
    16777215,   # 2  lt1      #FFFFFF
    6968388,    # 3  dk2      #44546A
    15132391,   # 4  lt2      #E7E6E6
    13998939,   # 5  accent1  #5B9BD5
    3243501,    # 6  accent2  #ED7D31
    10855845,   # 7  accent3  #A5A5A5
    49407,      # 8  accent4  #FFC000
    12874308,   # 9  accent5  #4472C4
    4697456,    # 10 accent6  #70AD47
    12673797,   # 11 hlink    #0563C1
    7491477     # 12 folHlink #954F72
)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeRgb[$i - 1]
}
